# Append the new "Milk - 2%" (twog) line item as row 6 on the active sheet,
# matching the existing rows which store every value (including the
# numeric-looking Quantity / Cost Per / Total Cost columns) as plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 6

$ws.Cells.Item($row, 1).Value = "twog"
$ws.Cells.Item($row, 2).Value = "Milk - 2%"
# Leading apostrophes force these numeric-looking strings to be stored as
# text (matching SKU/Quantity/Cost/Total columns elsewhere in the sheet)
# instead of being auto-coerced into numbers.
$ws.Cells.Item($row, 3).Value = "'14"
$ws.Cells.Item($row, 4).Value = "'16.08"
$ws.Cells.Item($row, 5).Value = "'225.12"

# Reset to the default style so the new row doesn't pick up the implicit
# "quote prefix" formatting the apostrophe entry would otherwise apply.
$ws.Range("A6:E6").Style = "Normal"
